# Append a new daily profit row (row 42) to the sheet, matching the
# existing pattern of text-formatted dates in column A and numeric
# profit values in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Range("A42")

# Force the cell to be treated as text first so the date-like string
# "09/28/2025" is stored as a literal string (like the other date
# cells in the column) instead of being auto-converted to a date
# serial number.
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/28/2025"
# Drop the temporary text formatting again so the new cell ends up
# with no explicit style, consistent with the surrounding rows.
$dateCell.ClearFormats()

$ws.Range("B42").Value = 14139.33
